# Todd update syllabus and TAM
# Applies the syllabus content revisions: renumbered "Lecture NN:" topic
# labels, reshuffled day-2/3/4 schedule items, a couple of corrected
# time slots, a new Day-3 reading reference, bolded "Lecture" topic rows,
# and small view/formatting touches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Monday ----
$ws.Range("D3").Value  = "Lecture 01: What are models? Why do we use them? Intro to Netlogo"
$ws.Range("D5").Value  = "Lecture 02: Conceptual Modeling"

# ---- Tuesday ----
$ws.Range("D12").Value = "Design presentations (1 slide)"
$ws.Range("D13").Value = "Lecture 03: Software testing, Building an ABM"
$ws.Range("C14").Value = 1000
$ws.Range("D14").Value = "Exercise: building an ABM"
$ws.Range("D16").Value = "Lecture 04: Design concepts: Agent sets, emergence, sensing, interactions"
$ws.Range("D17").Value = "Lecture 05: Becoming a modeler (Ms. Carra Carrillo)"
$ws.Range("C18").Value = 1430
$ws.Range("D18").Value = "Project work: Implementing your own ABM"

# ---- Wednesday ----
$ws.Range("D22").Value = "Lecture 06: Adaptive behavior, scheduling, collectives"
$ws.Range("E23").Value = "Ch. 16"
$ws.Range("D25").Value = "Lecture 07: Model Evaluation & pattern oriented modeling"

# ---- Thursday ----
$ws.Range("D31").Value = "Lecture 08:Uncertainty (Dr. Adam Duarte)"
$ws.Range("D34").Value = "Lecture 09: Communicating and documenting models"

# ---- Friday ----
$ws.Range("D37").Value = "Break for evening (Homework: prepare presentation)"

# ---- Bold the "Lecture NN" topic rows and give them the taller header-style
#      row height, matching the other lecture-topic row already in the sheet ----
foreach ($r in 16, 22, 25, 31, 34) {
    $ws.Range("D$r").Font.Bold = $true
    $ws.Rows("$r").RowHeight = 15.75
}

# ---- Column E needed to widen slightly to fit the new reading references ----
$ws.Columns("E").ColumnWidth = 24.75

# ---- Restore the selection to where the edits were focused ----
$ws.Range("D27").Select()
